# PH223_Lecture_17.pptx edit: "almost done with F2024"
#
# The change deletes the single picture on the last slide (sldId 261,
# which is slide index 34 in the deck) — shape id 2 / "Picture 2".
# (Date placeholder fields and revisionInfo timestamps are recalculated
# automatically by the runtime from its pinned clock, so no manual edit
# of those is required here.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(34)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Type -eq 13) {  # msoPicture
        $shape.Delete()
    }
}
